$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (24 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1618.7778
$ws.Range("I4").Value = 847.8
$ws.Range("K4").Value = 847.8
$ws.Range("M4").Value = -733.8
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H100").Value = 3352.5
$ws.Range("I100").Value = 1705
$ws.Range("K100").Value = 1705
$ws.Range("M100").Value = -1164
$ws.Range("H127").Value = 2021.25
$ws.Range("I127").Value = 1350.7142
$ws.Range("K127").Value = 4052.1426
$ws.Range("M127").Value = 907.8574000000003
$ws.Range("H137").Value = 3499
$ws.Range("J137").Value = 3499
$ws.Range("L137").Value = 10497
$ws.Range("N137").Value = -15597

# --- Sheet: ARM (46 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2653.5
$ws.Range("I45").Value = 2496.889
$ws.Range("K45").Value = 2496.889
$ws.Range("M45").Value = -2119.889
$ws.Range("H61").Value = 11497457
$ws.Range("I61").Value = 19610456
$ws.Range("J61").Value = 4041.5
$ws.Range("K61").Value = 19610456
$ws.Range("L61").Value = 4041.5
$ws.Range("M61").Value = -19610244
$ws.Range("N61").Value = -4465.5
$ws.Range("H74").Value = 2517.65
$ws.Range("I74").Value = 2484.2
$ws.Range("J74").Value = 2618
$ws.Range("K74").Value = 2484.2
$ws.Range("L74").Value = 2618
$ws.Range("M74").Value = -1610.2
$ws.Range("N74").Value = -4366
$ws.Range("H77").Value = 2517.65
$ws.Range("I77").Value = 2484.2
$ws.Range("J77").Value = 2618
$ws.Range("K77").Value = 12421
$ws.Range("L77").Value = 13090
$ws.Range("M77").Value = -8053
$ws.Range("N77").Value = -21826
$ws.Range("H110").Value = 200000800
$ws.Range("I110").Value = 333333820
$ws.Range("K110").Value = 333333820
$ws.Range("M110").Value = -333331775
$ws.Range("I132").Value = 34484976
$ws.Range("J132").Value = 1799
$ws.Range("K132").Value = 103454928
$ws.Range("L132").Value = 5397
$ws.Range("M132").Value = -103452398
$ws.Range("N132").Value = -10457
$ws.Range("H136").Value = 11497457
$ws.Range("I136").Value = 19610456
$ws.Range("J136").Value = 4041.5
$ws.Range("K136").Value = 58831368
$ws.Range("L136").Value = 12124.5
$ws.Range("M136").Value = -58828818
$ws.Range("N136").Value = -17224.5
$ws.Range("H139").Value = 90999.14
$ws.Range("J139").Value = 90999.14
$ws.Range("L139").Value = 90999.14
$ws.Range("N139").Value = -101279.14

# --- Sheet: BSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1404.5
$ws.Range("I12").Value = 478.33334
$ws.Range("J12").Value = 1960.2
$ws.Range("K12").Value = 478.33334
$ws.Range("L12").Value = 1960.2
$ws.Range("M12").Value = -310.33334
$ws.Range("N12").Value = -2296.2
$ws.Range("H37").Value = 371.77777
$ws.Range("I37").Value = 355.75
$ws.Range("K37").Value = 355.75
$ws.Range("M37").Value = -218.75
$ws.Range("H99").Value = 911.2
$ws.Range("I99").Value = 786.5
$ws.Range("K99").Value = 786.5
$ws.Range("M99").Value = 711.5

# --- Sheet: CRP (18 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2901.6572
$ws.Range("I31").Value = 1489.0476
$ws.Range("J31").Value = 3507.0613
$ws.Range("K31").Value = 1489.0476
$ws.Range("L31").Value = 3507.0613
$ws.Range("M31").Value = -1194.0476
$ws.Range("N31").Value = -4097.061299999999
$ws.Range("H34").Value = 2901.6572
$ws.Range("I34").Value = 1489.0476
$ws.Range("J34").Value = 3507.0613
$ws.Range("K34").Value = 1489.0476
$ws.Range("L34").Value = 3507.0613
$ws.Range("M34").Value = -1287.0476
$ws.Range("N34").Value = -3911.0613
$ws.Range("H139").Value = 77000
$ws.Range("J139").Value = 77000
$ws.Range("L139").Value = 77000
$ws.Range("N139").Value = -87280

# --- Sheet: CUL (30 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 287.5
$ws.Range("I26").Value = 75
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 225
$ws.Range("L26").Value = 1500
$ws.Range("M26").Value = 63
$ws.Range("N26").Value = -2076
$ws.Range("H99").Value = 2268.75
$ws.Range("I99").Value = 2268.75
$ws.Range("K99").Value = 6806.25
$ws.Range("M99").Value = -4560.25
$ws.Range("H129").Value = 2220.55
$ws.Range("I129").Value = 1125.5714
$ws.Range("K129").Value = 3376.7142
$ws.Range("M129").Value = 1623.2858
$ws.Range("H134").Value = 2638.1667
$ws.Range("I134").Value = 2638.1667
$ws.Range("K134").Value = 7914.500100000001
$ws.Range("M134").Value = -2844.500100000001
$ws.Range("H139").Value = 3041167.2
$ws.Range("I139").Value = 3716148.8
$ws.Range("K139").Value = 11148446.4
$ws.Range("M139").Value = -11143306.4
$ws.Range("H140").Value = 1585.0526
$ws.Range("I140").Value = 838.8333
$ws.Range("J140").Value = 2864.2856
$ws.Range("K140").Value = 2516.4999
$ws.Range("L140").Value = 8592.856800000001
$ws.Range("M140").Value = 2663.5001
$ws.Range("N140").Value = -18952.8568

# --- Sheet: GSM (29 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 7771.143
$ws.Range("I41").Value = 5849.75
$ws.Range("K41").Value = 5849.75
$ws.Range("M41").Value = -5494.75
$ws.Range("H63").Value = 39001
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 39001
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 39001
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -40373
$ws.Range("H66").Value = 39001
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 39001
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 117003
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -123867
$ws.Range("H113").Value = 3010.7646
$ws.Range("I113").Value = 1400
$ws.Range("K113").Value = 1400
$ws.Range("M113").Value = 770
$ws.Range("H132").Value = 3636.3225
$ws.Range("I132").Value = 3036.65
$ws.Range("J132").Value = 4726.636
$ws.Range("K132").Value = 9109.950000000001
$ws.Range("L132").Value = 14179.908
$ws.Range("M132").Value = -6579.950000000001
$ws.Range("N132").Value = -19239.908

# --- Sheet: LTW (65 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12823429
$ws.Range("I7").Value = 17243818
$ws.Range("J7").Value = 4303.9
$ws.Range("K7").Value = 17243818
$ws.Range("L7").Value = 4303.9
$ws.Range("M7").Value = -17243706
$ws.Range("N7").Value = -4527.9
$ws.Range("H9").Value = 1216
$ws.Range("I9").Value = 686
$ws.Range("J9").Value = 2011
$ws.Range("K9").Value = 686
$ws.Range("L9").Value = 2011
$ws.Range("M9").Value = -462
$ws.Range("N9").Value = -2459
$ws.Range("H13").Value = 3007
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3007
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3007
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3287
$ws.Range("H16").Value = 480.46155
$ws.Range("I16").Value = 480.46155
$ws.Range("K16").Value = 480.46155
$ws.Range("M16").Value = -310.46155
$ws.Range("H19").Value = 2088.5
$ws.Range("I19").Value = 995
$ws.Range("K19").Value = 995
$ws.Range("M19").Value = -825
$ws.Range("H31").Value = 698.6667
$ws.Range("I31").Value = 761.625
$ws.Range("J31").Value = 195
$ws.Range("K31").Value = 761.625
$ws.Range("L31").Value = 195
$ws.Range("M31").Value = -513.625
$ws.Range("N31").Value = -691
$ws.Range("H32").Value = 3063
$ws.Range("I32").Value = 3063
$ws.Range("K32").Value = 3063
$ws.Range("M32").Value = -2746
$ws.Range("H34").Value = 9761.25
$ws.Range("I34").Value = 9673.666999999999
$ws.Range("J34").Value = 10024
$ws.Range("K34").Value = 9673.666999999999
$ws.Range("L34").Value = 10024
$ws.Range("M34").Value = -9501.666999999999
$ws.Range("N34").Value = -10368
$ws.Range("H126").Value = 12823429
$ws.Range("I126").Value = 17243818
$ws.Range("J126").Value = 4303.9
$ws.Range("K126").Value = 51731454
$ws.Range("L126").Value = 12911.7
$ws.Range("M126").Value = -51728984
$ws.Range("N126").Value = -17851.7
$ws.Range("H132").Value = 5451.706
$ws.Range("I132").Value = 3585.8125
$ws.Range("J132").Value = 7110.278
$ws.Range("K132").Value = 10757.4375
$ws.Range("L132").Value = 21330.834
$ws.Range("M132").Value = -8227.4375
$ws.Range("N132").Value = -26390.834
$ws.Range("H136").Value = 2348.5293
$ws.Range("I136").Value = 1939.0741
$ws.Range("K136").Value = 5817.2223
$ws.Range("M136").Value = -3267.2223

# --- Sheet: WVR (43 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 13437.5
$ws.Range("I22").Value = 3500
$ws.Range("J22").Value = 23375
$ws.Range("K22").Value = 3500
$ws.Range("L22").Value = 23375
$ws.Range("M22").Value = -3207
$ws.Range("N22").Value = -23961
$ws.Range("H23").Value = 2575.2144
$ws.Range("I23").Value = 1988.8
$ws.Range("J23").Value = 4041.25
$ws.Range("K23").Value = 1988.8
$ws.Range("L23").Value = 4041.25
$ws.Range("M23").Value = -1759.8
$ws.Range("N23").Value = -4499.25
$ws.Range("H64").Value = 35034
$ws.Range("I64").Value = 25051
$ws.Range("J64").Value = 55000
$ws.Range("K64").Value = 25051
$ws.Range("L64").Value = 55000
$ws.Range("M64").Value = -24803
$ws.Range("N64").Value = -55496
$ws.Range("H67").Value = 35034
$ws.Range("I67").Value = 25051
$ws.Range("J67").Value = 55000
$ws.Range("K67").Value = 25051
$ws.Range("L67").Value = 55000
$ws.Range("M67").Value = -24193
$ws.Range("N67").Value = -56716
$ws.Range("H107").Value = 544.8333
$ws.Range("I107").Value = 554.5333000000001
$ws.Range("K107").Value = 1663.5999
$ws.Range("M107").Value = 256.4000999999998
$ws.Range("H122").Value = 1865
$ws.Range("I122").Value = 1821.7142
$ws.Range("J122").Value = 2067
$ws.Range("K122").Value = 5465.142599999999
$ws.Range("L122").Value = 6201
$ws.Range("M122").Value = -3015.142599999999
$ws.Range("N122").Value = -11101
$ws.Range("H132").Value = 5388.3335
$ws.Range("I132").Value = 5185.3335
$ws.Range("K132").Value = 15556.0005
$ws.Range("M132").Value = -13026.0005

Write-Output "Applied 270 cell updates"